# Updated run for publication
# Applies updated frequency values to the "LTR-E" frequency table on Sheet1.
# Source diff recomputed the per-position base frequencies using a new
# alignment denominator (430 sequences instead of 77), so most non-zero /
# near-zero cells in rows 2-5 (B:X) receive refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.00930232558139535
$ws.Range("D2").Value = 0.0186046511627907
$ws.Range("F2").Value = 0.00697674418604651
$ws.Range("G2").Value = 0.00465116279069767
$ws.Range("H2").Value = 0.351162790697674
$ws.Range("I2").Value = 0.00697674418604651
$ws.Range("K2").Value = 0.0116279069767442
$ws.Range("L2").Value = 0.993023255813953
$ws.Range("N2").Value = 0.193023255813953
$ws.Range("Q2").Value = 0.995348837209302
$ws.Range("S2").Value = 0.913953488372093
$ws.Range("T2").Value = 0.00232558139534884
$ws.Range("U2").Value = 0.988372093023256
$ws.Range("V2").Value = 0.92093023255814
$ws.Range("W2").Value = 0.00232558139534884
$ws.Range("X2").Value = 0.00697674418604651
$ws.Range("B3").Value = 0.809302325581395
$ws.Range("D3").Value = 0.00232558139534884
$ws.Range("E3").Value = 0.00697674418604651
$ws.Range("H3").Value = 0.00930232558139535
$ws.Range("I3").Value = 0.00930232558139535
$ws.Range("J3").Value = 0.909302325581395
$ws.Range("K3").Value = 0.00465116279069767
$ws.Range("L3").Value = 0.00232558139534884
$ws.Range("M3").Value = 0.92093023255814
$ws.Range("N3").Value = 0.8
$ws.Range("P3").Value = 0.927906976744186
$ws.Range("R3").Value = 0.997674418604651
$ws.Range("S3").Value = 0.00232558139534884
$ws.Range("T3").Value = 0.995348837209302
$ws.Range("U3").Value = 0.00697674418604651
$ws.Range("V3").Value = 0.00232558139534884
$ws.Range("B4").Value = 0.00930232558139535
$ws.Range("C4").Value = 0.00465116279069767
$ws.Range("D4").Value = 0.965116279069767
$ws.Range("E4").Value = 0.00465116279069767
$ws.Range("F4").Value = 0.993023255813953
$ws.Range("G4").Value = 0.995348837209302
$ws.Range("H4").Value = 0.637209302325581
$ws.Range("J4").Value = 0.0116279069767442
$ws.Range("L4").Value = 0.00465116279069767
$ws.Range("P4").Value = 0.00465116279069767
$ws.Range("Q4").Value = 0.00232558139534884
$ws.Range("R4").Value = 0.00232558139534884
$ws.Range("S4").Value = 0.0837209302325581
$ws.Range("T4").Value = 0.00232558139534884
$ws.Range("U4").Value = 0.00465116279069767
$ws.Range("V4").Value = 0.0767441860465116
$ws.Range("W4").Value = 0.995348837209302
$ws.Range("X4").Value = 0.993023255813953
$ws.Range("B5").Value = 0.181395348837209
$ws.Range("C5").Value = 0.986046511627907
$ws.Range("D5").Value = 0.013953488372093
$ws.Range("E5").Value = 0.988372093023256
$ws.Range("H5").Value = 0.00232558139534884
$ws.Range("I5").Value = 0.983720930232558
$ws.Range("J5").Value = 0.0790697674418605
$ws.Range("K5").Value = 0.983720930232558
$ws.Range("M5").Value = 0.0767441860465116
$ws.Range("N5").Value = 0.00697674418604651
$ws.Range("P5").Value = 0.0674418604651163
$ws.Range("Q5").Value = 0.00232558139534884
$ws.Range("W5").Value = 0.00232558139534884
